# Trading update: 2026-02-18 10:40:04
#
# New MarketMaking trades landed (#32-#36). The previously "latest" open
# trade (#31, row 32 on "All Trades") loses its extended live-tracking
# metadata (Capital After / Slippage / Confidence / Entry Reason /
# Duration) and gets an Exit Price of 0 recorded, while the new trades
# take over that extended-metadata role. The "MarketMaking" strategy
# sheet is refreshed to mirror only the newest trade plus the rest of
# the incoming batch.

$wb = $excel.ActiveWorkbook

$allTrades = $wb.Worksheets.Item("All Trades")
$marketMaking = $wb.Worksheets.Item("MarketMaking")

function Set-TradeRow($ws, $row, $tradeNum, $date, $time, $strategy, $side, $entryPrice, $status, $plPct, $plDollar, $capitalAfter, $entrySlippage, $exitSlippage, $confidence, $entryReason, $duration) {
    $ws.Cells.Item($row, 1).Value = $tradeNum

    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $date

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $time

    $ws.Cells.Item($row, 4).Value = $strategy
    $ws.Cells.Item($row, 5).Value = $side
    $ws.Cells.Item($row, 6).Value = $entryPrice
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = $status
    $ws.Cells.Item($row, 9).Value = $plPct
    $ws.Cells.Item($row, 10).Value = $plDollar
    $ws.Cells.Item($row, 11).Value = $capitalAfter
    $ws.Cells.Item($row, 12).Value = $entrySlippage
    $ws.Cells.Item($row, 13).Value = $exitSlippage
    $ws.Cells.Item($row, 14).Value = $confidence
    $ws.Cells.Item($row, 15).Value = $entryReason
    $ws.Cells.Item($row, 16).Value = ""
    $ws.Cells.Item($row, 17).Value = $duration
}

# --- "All Trades": trade #31 (row 32) rolls off the extended-metadata
#     tracking slot: Exit Price becomes 0, and Capital After / Slippage /
#     Confidence / Entry Reason / Duration go blank again. ---
$allTrades.Cells.Item(32, 7).Value = 0
$allTrades.Cells.Item(32, 11).Value = ""
$allTrades.Cells.Item(32, 12).Value = ""
$allTrades.Cells.Item(32, 13).Value = ""
$allTrades.Cells.Item(32, 14).Value = ""
$allTrades.Cells.Item(32, 15).Value = ""
$allTrades.Cells.Item(32, 17).Value = ""

# --- New trades #32-#36 appended to "All Trades" (rows 33-37) ---
Set-TradeRow $allTrades 33 32 "2026-02-18" "10:39:05" "MarketMaking" "UP" 0.9 "OPEN" 0 0 100 0 0 0.6 "Normal spread capture: 714 bps" 0
Set-TradeRow $allTrades 34 33 "2026-02-18" "10:39:11" "MarketMaking" "DOWN" 0.1 "OPEN" 0 0 100 0 0 0.6 "Normal spread capture: 714 bps" 0
Set-TradeRow $allTrades 35 34 "2026-02-18" "10:39:17" "MarketMaking" "UP" 0.91 "OPEN" 0 0 100 0 0 0.6 "Normal spread capture: 714 bps" 0
Set-TradeRow $allTrades 36 35 "2026-02-18" "10:39:23" "MarketMaking" "DOWN" 0.08 "OPEN" 0 0 100 0 0 0.6 "Normal spread capture: 714 bps" 0
Set-TradeRow $allTrades 37 36 "2026-02-18" "10:39:30" "MarketMaking" "DOWN" 0.06 "OPEN" 0 0 100 0 0 0.6 "Normal spread capture: 714 bps" 0

# --- "MarketMaking" strategy sheet: refresh row 2 (was trade #31) to
#     trade #32's details, then append trades #33-#36 in rows 3-6. ---
$marketMaking.Cells.Item(2, 1).Value = 32
$marketMaking.Cells.Item(2, 3).NumberFormat = "@"
$marketMaking.Cells.Item(2, 3).Value = "10:39:05"
$marketMaking.Cells.Item(2, 6).Value = 0.9
$marketMaking.Cells.Item(2, 15).Value = "Normal spread capture: 714 bps"

Set-TradeRow $marketMaking 3 33 "2026-02-18" "10:39:11" "MarketMaking" "DOWN" 0.1 "OPEN" 0 0 100 0 0 0.6 "Normal spread capture: 714 bps" 0
Set-TradeRow $marketMaking 4 34 "2026-02-18" "10:39:17" "MarketMaking" "UP" 0.91 "OPEN" 0 0 100 0 0 0.6 "Normal spread capture: 714 bps" 0
Set-TradeRow $marketMaking 5 35 "2026-02-18" "10:39:23" "MarketMaking" "DOWN" 0.08 "OPEN" 0 0 100 0 0 0.6 "Normal spread capture: 714 bps" 0
Set-TradeRow $marketMaking 6 36 "2026-02-18" "10:39:30" "MarketMaking" "DOWN" 0.06 "OPEN" 0 0 100 0 0 0.6 "Normal spread capture: 714 bps" 0
